$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Append a new results row (row 8) to the Summary sheet, mirroring the
# structure of the existing rows (A..AF), for the new
# "GreenFieldHydro_Island_2040" simulation case.

$row = 8

$numericValues = @{
    "A" = 230085772776.3445
    "B" = 13313796935.41874
    "C" = 0
    "D" = 214741056080.9258
    "E" = 2030919760
    "F" = 228054853016.3445
    "G" = 2030919760
    "H" = 0
    "I" = 0
    "J" = 0
    "K" = 0
    "L" = 0
    "M" = 230085772776.3445
    "N" = 0
    "O" = 0
    "P" = 0
    "Q" = 0
    "R" = 0
    "S" = 0
    "T" = 0
    "U" = 0
    "V" = 1827.349999904633
    "W" = 230085772185.4904
    "X" = 230085772776.3445
    "Y" = 590.8540649414062
    "AB" = -1
    "AC" = -1
    "AD" = 1
}

foreach ($col in $numericValues.Keys) {
    $ws.Range("$col$row").Value = $numericValues[$col]
}

$textValues = @{
    "Z"  = "optimal"
    "AA" = "costs_emissionlimit"
    "AE" = "GreenFieldHydro_Island_2040"
    "AF" = "rawResults\20251124152705_GreenFieldHydro_Island_2040-1"
}

foreach ($col in $textValues.Keys) {
    $ws.Range("$col$row").Value = $textValues[$col]
}
